$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-blank "Value" (column D) entries for several game
# elements (rows 5-13). Column G mirrors these via a shared formula and
# the G2 total (SUM(G4:G23)) will recalculate from 50 to 100.
$ws.Range("D5").Value = 25
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 5
$ws.Range("D8").Value = 15
$ws.Range("D11").Value = 5
$ws.Range("D13").Value = 0

# Make sure all dependent formulas (G column, G2 total, etc.) are refreshed.
$excel.CalculateFull()

# Update the saved cursor/selection position on the sheet.
[void]$ws.Range("B11").Select()
